# Regenerate save_data: use K (Strike count) instead of Strike#,
# recompute std/mean-derived K values and write the resulting s_vals
# into column G (header 'K') for each trade row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values computed from the regenerated save_data (s_vals).
$kValues = @{
    2 = 0
    3 = 1
    4 = 1
    5 = 3
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 3
    11 = 1
    13 = 1
    14 = 0
    15 = 0
    18 = 1
    19 = 2
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 3
    39 = 1
    40 = 2
    41 = 0
    42 = 0
    43 = 3
    44 = 2
    45 = 1
    46 = 0
    47 = 4
    48 = 2
    49 = 2
    50 = 1
    51 = 1
    52 = 1
    54 = 0
    55 = 0
    56 = 1
    57 = 1
    58 = 3
    59 = 0
    60 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
